# Update the "metadata nomeclatures" template sheet: add a new header
# column G ("target element from (your domain standards)") and update
# the current selection/view to reflect where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata nomeclatures")

# Add the new header cell in G1, matching the style of the existing
# header row (bold/filled header, via the row's own formatting).
$ws.Range("G1").Value = "target element from (your domain standards)"

# Reflect the saved selection/scroll state left in the workbook.
$ws.Activate()
$ws.Range("F10").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
